$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2360.9395
$ws.Range("I70").Value = 3695.9412
$ws.Range("J70").Value = 942.5
$ws.Range("K70").Value = 11087.8236
$ws.Range("L70").Value = 2827.5
$ws.Range("M70").Value = -10817.8236
$ws.Range("N70").Value = -3367.5
$ws.Range("H73").Value = 2360.9395
$ws.Range("I73").Value = 3695.9412
$ws.Range("J73").Value = 942.5
$ws.Range("K73").Value = 11087.8236
$ws.Range("L73").Value = 2827.5
$ws.Range("M73").Value = -10151.8236
$ws.Range("N73").Value = -4699.5
$ws.Range("H86").Value = 4067.0667
$ws.Range("I86").Value = 1840.0555
$ws.Range("J86").Value = 7407.5835
$ws.Range("K86").Value = 1840.0555
$ws.Range("L86").Value = 7407.5835
$ws.Range("M86").Value = -717.0554999999999
$ws.Range("N86").Value = -9653.583500000001
$ws.Range("H89").Value = 4067.0667
$ws.Range("I89").Value = 1840.0555
$ws.Range("J89").Value = 7407.5835
$ws.Range("K89").Value = 9200.2775
$ws.Range("L89").Value = 37037.9175
$ws.Range("M89").Value = -3584.2775
$ws.Range("N89").Value = -48269.9175
$ws.Range("H127").Value = 622
$ws.Range("I127").Value = 271.81818
$ws.Range("J127").Value = 1050
$ws.Range("K127").Value = 815.45454
$ws.Range("L127").Value = 3150
$ws.Range("M127").Value = 4144.54546
$ws.Range("N127").Value = -13070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3012.6584
$ws.Range("I102").Value = 1548.0625
$ws.Range("J102").Value = 8220.111000000001
$ws.Range("K102").Value = 1548.0625
$ws.Range("L102").Value = 8220.111000000001
$ws.Range("M102").Value = 73.9375
$ws.Range("N102").Value = -11464.111
$ws.Range("H122").Value = 2153.9167
$ws.Range("I122").Value = 2186.6924
$ws.Range("K122").Value = 6560.0772
$ws.Range("M122").Value = -4110.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4841.385
$ws.Range("I99").Value = 5993.8
$ws.Range("K99").Value = 5993.8
$ws.Range("M99").Value = -4495.8
$ws.Range("H105").Value = 2044.7587
$ws.Range("I105").Value = 1944.3334
$ws.Range("J105").Value = 2209.0908
$ws.Range("K105").Value = 1944.3334
$ws.Range("L105").Value = 2209.0908
$ws.Range("M105").Value = -197.3334
$ws.Range("N105").Value = -5703.0908
$ws.Range("H107").Value = 1966.0435
$ws.Range("I107").Value = 1861.8334
$ws.Range("J107").Value = 2341.2
$ws.Range("K107").Value = 1861.8334
$ws.Range("L107").Value = 2341.2
$ws.Range("M107").Value = 58.16660000000002
$ws.Range("N107").Value = -6181.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 18052.9
$ws.Range("J28").Value = 18052.9
$ws.Range("L28").Value = 18052.9
$ws.Range("N28").Value = -18542.9
$ws.Range("H31").Value = 2504.0461
$ws.Range("I31").Value = 1533.375
$ws.Range("J31").Value = 4057.12
$ws.Range("K31").Value = 1533.375
$ws.Range("L31").Value = 4057.12
$ws.Range("M31").Value = -1238.375
$ws.Range("N31").Value = -4647.12
$ws.Range("H34").Value = 2504.0461
$ws.Range("I34").Value = 1533.375
$ws.Range("J34").Value = 4057.12
$ws.Range("K34").Value = 1533.375
$ws.Range("L34").Value = 4057.12
$ws.Range("M34").Value = -1331.375
$ws.Range("N34").Value = -4461.12
$ws.Range("H58").Value = 1556.575
$ws.Range("I58").Value = 967.6087
$ws.Range("J58").Value = 2353.4119
$ws.Range("K58").Value = 967.6087
$ws.Range("L58").Value = 2353.4119
$ws.Range("M58").Value = -764.6087
$ws.Range("N58").Value = -2759.4119
$ws.Range("H86").Value = 8008.0835
$ws.Range("I86").Value = 6371
$ws.Range("J86").Value = 10300
$ws.Range("K86").Value = 6371
$ws.Range("L86").Value = 10300
$ws.Range("M86").Value = -5248
$ws.Range("N86").Value = -12546
$ws.Range("H89").Value = 8008.0835
$ws.Range("I89").Value = 6371
$ws.Range("J89").Value = 10300
$ws.Range("K89").Value = 31855
$ws.Range("L89").Value = 51500
$ws.Range("M89").Value = -26239
$ws.Range("N89").Value = -62732
$ws.Range("H96").Value = 14078
$ws.Range("J96").Value = 14078
$ws.Range("L96").Value = 14078
$ws.Range("N96").Value = -19570
$ws.Range("H107").Value = 369.9
$ws.Range("I107").Value = 285.7143
$ws.Range("J107").Value = 566.3333
$ws.Range("K107").Value = 285.7143
$ws.Range("L107").Value = 566.3333
$ws.Range("M107").Value = 1634.2857
$ws.Range("N107").Value = -4406.3333
$ws.Range("H136").Value = 1556.575
$ws.Range("I136").Value = 967.6087
$ws.Range("J136").Value = 2353.4119
$ws.Range("K136").Value = 2902.8261
$ws.Range("L136").Value = 7060.2357
$ws.Range("M136").Value = -352.8261000000002
$ws.Range("N136").Value = -12160.2357

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 872.95
$ws.Range("I5").Value = 644.875
$ws.Range("J5").Value = 1025
$ws.Range("K5").Value = 1934.625
$ws.Range("L5").Value = 3075
$ws.Range("M5").Value = -1822.625
$ws.Range("N5").Value = -3299
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H129").Value = 1707.3
$ws.Range("I129").Value = 1682.2222
$ws.Range("J129").Value = 1933
$ws.Range("K129").Value = 5046.6666
$ws.Range("L129").Value = 5799
$ws.Range("M129").Value = -46.66659999999956
$ws.Range("N129").Value = -15799
$ws.Range("H135").Value = 872.95
$ws.Range("I135").Value = 644.875
$ws.Range("J135").Value = 1025
$ws.Range("K135").Value = 5803.875
$ws.Range("L135").Value = 9225
$ws.Range("M135").Value = -3268.875
$ws.Range("N135").Value = -14295
$ws.Range("H138").Value = 2141.7222
$ws.Range("I138").Value = 1801.8182
$ws.Range("K138").Value = 5405.4546
$ws.Range("M138").Value = -265.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3368
$ws.Range("I132").Value = 4141.619
$ws.Range("J132").Value = 2766.2964
$ws.Range("K132").Value = 12424.857
$ws.Range("L132").Value = 8298.889200000001
$ws.Range("M132").Value = -9894.857
$ws.Range("N132").Value = -13358.8892

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2673.4375
$ws.Range("I7").Value = 2327.8572
$ws.Range("K7").Value = 2327.8572
$ws.Range("M7").Value = -2215.8572
$ws.Range("H16").Value = 1990.5385
$ws.Range("I16").Value = 1978.7142
$ws.Range("J16").Value = 2040.2
$ws.Range("K16").Value = 1978.7142
$ws.Range("L16").Value = 2040.2
$ws.Range("M16").Value = -1808.7142
$ws.Range("N16").Value = -2380.2
$ws.Range("H40").Value = 2685.4211
$ws.Range("I40").Value = 2558.7856
$ws.Range("J40").Value = 3040
$ws.Range("K40").Value = 2558.7856
$ws.Range("L40").Value = 3040
$ws.Range("M40").Value = -2422.7856
$ws.Range("N40").Value = -3312
$ws.Range("H122").Value = 3349.8
$ws.Range("I122").Value = 2820.8
$ws.Range("J122").Value = 3614.3
$ws.Range("K122").Value = 8462.400000000001
$ws.Range("L122").Value = 10842.9
$ws.Range("M122").Value = -6012.400000000001
$ws.Range("N122").Value = -15742.9
$ws.Range("H126").Value = 2673.4375
$ws.Range("I126").Value = 2327.8572
$ws.Range("K126").Value = 6983.571599999999
$ws.Range("M126").Value = -4513.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1940
$ws.Range("I81").Value = 1426.1305
$ws.Range("J81").Value = 3628.4285
$ws.Range("K81").Value = 2852.261
$ws.Range("L81").Value = 7256.857
$ws.Range("M81").Value = -1791.261
$ws.Range("N81").Value = -9378.857
$ws.Range("H84").Value = 1940
$ws.Range("I84").Value = 1426.1305
$ws.Range("J84").Value = 3628.4285
$ws.Range("K84").Value = 14261.305
$ws.Range("L84").Value = 36284.285
$ws.Range("M84").Value = -8957.305
$ws.Range("N84").Value = -46892.285
$ws.Range("H126").Value = 1046.9032
$ws.Range("I126").Value = 833.76
$ws.Range("J126").Value = 1935
$ws.Range("K126").Value = 2501.28
$ws.Range("L126").Value = 5805
$ws.Range("M126").Value = -31.27999999999975
$ws.Range("N126").Value = -10745
$ws.Range("H136").Value = 25662480
$ws.Range("I136").Value = 34518960
$ws.Range("J136").Value = 14495614
$ws.Range("K136").Value = 103556880
$ws.Range("L136").Value = 43486842
$ws.Range("M136").Value = -103554330
$ws.Range("N136").Value = -43491942
